# Refresh the cryptocurrency price / 1h-volume figures, and swap the
# Monero/Aptos rows (Aptos moves above Monero in the ranking).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.688.04"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "2.616.12"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("E4").Value = "  -0.65%  "
$ws.Range("D5").Value = "'516.03"
$ws.Range("E5").Value = "  +1.33%  "
$ws.Range("D6").Value = "'154.66"
$ws.Range("E6").Value = "  -1.23%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "'0.597"
$ws.Range("E8").Value = "  +1.48%  "
$ws.Range("D9").Value = "2.629.88"
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("D10").Value = "'6.69"
$ws.Range("E10").Value = "  +3.96%  "
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("E13").Value = "  +1.82%  "
$ws.Range("D14").Value = "3.074.60"
$ws.Range("E14").Value = "  -0.58%  "
$ws.Range("D15").Value = "60.727.02"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").Value = "'21.77"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("D18").Value = "2.625.34"
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("D19").Value = "'4.75"
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("D20").Value = "'357.76"
$ws.Range("E20").Value = "  +3.84%  "
$ws.Range("D21").Value = "'10.66"
$ws.Range("E21").Value = "  +2.09%  "
$ws.Range("D22").Value = "'6.23"
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "'61.18"
$ws.Range("E24").Value = "  +1.25%  "
$ws.Range("D25").Value = "'0.428"
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("D26").Value = "2.736.80"
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("D28").Value = "'0.996"
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("D29").Value = "0.0₃0848"
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("D30").Value = "'7.38"
$ws.Range("E30").Value = "  -2.17%  "
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").Value = "'19.46"
$ws.Range("E32").Value = "  +0.47%  "
$ws.Range("E33").Value = "  +0.85%  "
$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").Value = "'5.92"
$ws.Range("E34").Value = "  +3.56%  "
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").Value = "'150.88"
$ws.Range("E35").Value = "  -3.50%  "
$ws.Range("D36").Value = "'4.04"
$ws.Range("E36").Value = "  +0.84%  "
$ws.Range("D37").Value = "'1.19"
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("D38").Value = "'0.891"
$ws.Range("E38").Value = "  +6.10%  "
$ws.Range("E39").Value = "  +1.04%  "
$ws.Range("D40").Value = "'0.852"
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("D41").Value = "'36.43"
$ws.Range("E41").Value = "  +2.11%  "
$ws.Range("E42").Value = "  -1.12%  "
$ws.Range("D43").Value = "'291.10"
$ws.Range("E43").Value = "  -4.44%  "
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("D46").Value = "'0.0558"
$ws.Range("E46").Value = "  -2.10%  "
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("D48").Value = "'19.84"
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("E49").Value = "  +2.37%  "
$ws.Range("E50").Value = "  +0.49%  "
$ws.Range("E51").Value = "  +0.07%  "
